$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Fix capitalization of "de"/"la" -> "De"/"La" in municipality names ---
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("B6").Value = "Mazapa De Madero"
$ws.Range("B15").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B22").Value = "Tecali De Herrera"
$ws.Range("B35").Value = "Martínez De La Torre"

# --- Remove trailing metadata/footer rows (44-48 and 476-480) ---
$ws.Rows("476:480").Delete()
$ws.Rows("44:48").Delete()
